# Update localization status report: mark two files as "In Translation"
# instead of "Ready for handoff" (files 49145225-...md and b1b21d64-...md)
# across the Overview sheet and the per-locale (zh-cn, de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale status
# Row 3 -> 49145225-09c9-46ae-803f-739b042b9bda.md
# Row 4 -> b1b21d64-4f4b-4615-91cd-ab12a367e019.md
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

# zh-cn sheet: column C holds the Status field
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

# de-de sheet: column C holds the Status field
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
